$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the row above down into the new row, then overwrite values
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)   # xlPasteFormats

# Add the new changelog entry row (44295 = 2021-04-09 as an Excel date serial)
$ws.Range("A6").Value2 = 44295
$ws.Range("B6").Value2 = "1.0.0"
$ws.Range("C6").Value2 = "First official release of DefiChain-Analytics"

# Update the active selection, mirroring the diff's recorded selection state
$ws.Range("C9").Select()
